$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 272
$ws1.Range("F5").Value = 3339
$ws1.Range("F6").Value = 2129
$ws1.Range("F9").Value = 38
$ws1.Range("F11").Value = 1234
$ws1.Range("F12").Value = 222
$ws1.Range("F13").Value = 1347

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 272
$ws4.Range("F5").Value = 3339
$ws4.Range("F6").Value = 2129
$ws4.Range("F10").Value = 38
$ws4.Range("F14").Value = 1235
$ws4.Range("F15").Value = 222
$ws4.Range("F16").Value = 1347
